$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.411.66"
$ws.Range("E2").Value = "  +4.05%  "
$ws.Range("D3").Value = "3.180.30"
$ws.Range("E3").Value = "  +4.76%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.241"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +17.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.71%  "
$ws.Range("D10").Value = "3.187.10"
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.602"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +37.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000255"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +32.95%  "
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "3.771.51"
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +11.82%  "
$ws.Range("D17").Value = "79.441.18"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "3.194.25"
$ws.Range("E18").Value = "  +5.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +26.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +20.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.81%  "
$ws.Range("D25").Value = "3.350.07"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "77.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000123"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "542.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.152"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +30.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.121"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.98%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.411"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "192.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.21%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.799"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +17.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.638"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.55%  "
